$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$rng = $footer.Range.Duplicate
$found = $rng.Find.Execute("gatan 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
Write-Host "Found: $found"
Write-Host "RngText: $($rng.Text)"
Write-Host "Start: $($rng.Start) End: $($rng.End)"
$rng.InsertAfter(", 118 27")
Write-Host $footer.Range.Text
